$wb = $excel.ActiveWorkbook

# Sheet ALC, row 9
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 725
$ws.Range("I9").Value = 736.8421
$ws.Range("J9").Value = 500
$ws.Range("K9").Value = 736.8421
$ws.Range("L9").Value = 500
$ws.Range("M9").Value = -567.8421
$ws.Range("N9").Value = -838

# Sheet ALC, row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4634311.5
$ws.Range("J17").Value = 5213488
$ws.Range("L17").Value = 15640464
$ws.Range("N17").Value = -15640800

# Sheet ALC, row 53
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 319.33334
$ws.Range("I53").Value = 392.66666
$ws.Range("J53").Value = 172.66667
$ws.Range("K53").Value = 392.66666
$ws.Range("L53").Value = 172.66667
$ws.Range("M53").Value = 244.33334
$ws.Range("N53").Value = -1446.66667

# Sheet ALC, row 92
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 537.3333
$ws.Range("I92").Value = 447.42856
$ws.Range("K92").Value = 447.42856
$ws.Range("M92").Value = 800.5714399999999

# Sheet ALC, row 96
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 993.6316
$ws.Range("I96").Value = 684.875
$ws.Range("J96").Value = 1218.1818
$ws.Range("K96").Value = 2054.625
$ws.Range("L96").Value = 3654.5454
$ws.Range("M96").Value = -681.625
$ws.Range("N96").Value = -6400.5454

# Sheet ALC, row 99
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 6445.5
$ws.Range("I99").Value = 1877.3334
$ws.Range("J99").Value = 20150
$ws.Range("K99").Value = 5632.0002
$ws.Range("L99").Value = 60450
$ws.Range("M99").Value = -4134.0002
$ws.Range("N99").Value = -63446

# Sheet ALC, row 115
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H115").Value = 626.46155
$ws.Range("I115").Value = 364.4
$ws.Range("K115").Value = 1093.2
$ws.Range("M115").Value = 473.8000000000002

# Sheet ALC, row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2086238.6
$ws.Range("I116").Value = 2607004.5
$ws.Range("J116").Value = 3174.75
$ws.Range("K116").Value = 2607004.5
$ws.Range("L116").Value = 3174.75
$ws.Range("M116").Value = -2603562.5
$ws.Range("N116").Value = -10058.75

# Sheet ALC, row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1691.0555
$ws.Range("I135").Value = 1200.1875
$ws.Range("K135").Value = 10801.6875
$ws.Range("M135").Value = -8266.6875

# Sheet ARM, row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1106.619
$ws.Range("I2").Value = 859.9286
$ws.Range("J2").Value = 1600
$ws.Range("K2").Value = 859.9286
$ws.Range("L2").Value = 1600
$ws.Range("M2").Value = -746.9286
$ws.Range("N2").Value = -1826

# Sheet ARM, row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1753.1875
$ws.Range("I45").Value = 1923
$ws.Range("J45").Value = 1379.6
$ws.Range("K45").Value = 1923
$ws.Range("L45").Value = 1379.6
$ws.Range("M45").Value = -1546
$ws.Range("N45").Value = -2133.6

# Sheet ARM, row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1779.2
$ws.Range("I102").Value = 1757.4166
$ws.Range("J102").Value = 1866.3334
$ws.Range("K102").Value = 1757.4166
$ws.Range("L102").Value = 1866.3334
$ws.Range("M102").Value = -135.4166
$ws.Range("N102").Value = -5110.3334

# Sheet ARM, row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1106.619
$ws.Range("I116").Value = 859.9286
$ws.Range("J116").Value = 1600
$ws.Range("K116").Value = 859.9286
$ws.Range("L116").Value = 1600
$ws.Range("M116").Value = 1434.0714
$ws.Range("N116").Value = -6188

# Sheet ARM, row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1902
$ws.Range("I122").Value = 1709.909
$ws.Range("J122").Value = 2465.4666
$ws.Range("K122").Value = 5129.727000000001
$ws.Range("L122").Value = 7396.399800000001
$ws.Range("M122").Value = -2679.727000000001
$ws.Range("N122").Value = -12296.3998

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 20448.166
$ws.Range("I132").Value = 32895.03
$ws.Range("K132").Value = 98685.09
$ws.Range("M132").Value = -96155.09

# Sheet BSM, row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1106.619
$ws.Range("I3").Value = 859.9286
$ws.Range("J3").Value = 1600
$ws.Range("K3").Value = 859.9286
$ws.Range("L3").Value = 1600
$ws.Range("M3").Value = -745.9286
$ws.Range("N3").Value = -1828

# Sheet BSM, row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 15522.81
$ws.Range("I94").Value = 10104.272
$ws.Range("J94").Value = 21483.2
$ws.Range("K94").Value = 10104.272
$ws.Range("L94").Value = 21483.2
$ws.Range("M94").Value = -9653.272000000001
$ws.Range("N94").Value = -22385.2

# Sheet BSM, row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2313.8125
$ws.Range("I99").Value = 2249.8333
$ws.Range("J99").Value = 2352.2
$ws.Range("K99").Value = 2249.8333
$ws.Range("L99").Value = 2352.2
$ws.Range("M99").Value = -751.8332999999998
$ws.Range("N99").Value = -5348.2

# Sheet BSM, row 103
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 18328.5
$ws.Range("J103").Value = 18328.5
$ws.Range("L103").Value = 18328.5
$ws.Range("N103").Value = -20672.5

# Sheet BSM, row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2093.5
$ws.Range("I105").Value = 1805.8334
$ws.Range("J105").Value = 2525
$ws.Range("K105").Value = 1805.8334
$ws.Range("L105").Value = 2525
$ws.Range("M105").Value = -58.83339999999998
$ws.Range("N105").Value = -6019

# Sheet BSM, row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 257589.83
$ws.Range("I134").Value = 313290.78
$ws.Range("J134").Value = 2956.8572
$ws.Range("K134").Value = 939872.3400000001
$ws.Range("L134").Value = 8870.571599999999
$ws.Range("M134").Value = -937337.3400000001
$ws.Range("N134").Value = -13940.5716

# Sheet CUL, row 12
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 51.074074
$ws.Range("I12").Value = 82.28570999999999
$ws.Range("J12").Value = 40.15
$ws.Range("K12").Value = 246.85713
$ws.Range("L12").Value = 120.45
$ws.Range("M12").Value = -73.85712999999998
$ws.Range("N12").Value = -466.45

# Sheet CUL, row 38
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 111.666664
$ws.Range("I38").Value = 133.5
$ws.Range("J38").Value = 68
$ws.Range("K38").Value = 400.5
$ws.Range("L38").Value = 204
$ws.Range("M38").Value = -53.5
$ws.Range("N38").Value = -898

# Sheet CUL, row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 221.28572
$ws.Range("I107").Value = 200
$ws.Range("J107").Value = 237.25
$ws.Range("K107").Value = 600
$ws.Range("L107").Value = 711.75
$ws.Range("M107").Value = 1320
$ws.Range("N107").Value = -4551.75

# Sheet CUL, row 109
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 2063.7104
$ws.Range("I109").Value = 919.2273
$ws.Range("J109").Value = 3637.375
$ws.Range("K109").Value = 2757.6819
$ws.Range("L109").Value = 10912.125
$ws.Range("M109").Value = -1717.6819
$ws.Range("N109").Value = -12992.125

# Sheet CUL, row 129
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 5556383.5
$ws.Range("I129").Value = 331.5625
$ws.Range("J129").Value = 11906158
$ws.Range("K129").Value = 994.6875
$ws.Range("L129").Value = 35718474
$ws.Range("M129").Value = 4005.3125
$ws.Range("N129").Value = -35728474

# Sheet CUL, row 134
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 968.625
$ws.Range("I134").Value = 731.9286
$ws.Range("J134").Value = 1300
$ws.Range("K134").Value = 2195.7858
$ws.Range("L134").Value = 3900
$ws.Range("M134").Value = 2874.2142
$ws.Range("N134").Value = -14040

# Sheet GSM, row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1300.8462
$ws.Range("I122").Value = 1156.4482
$ws.Range("J122").Value = 1719.6
$ws.Range("K122").Value = 3469.3446
$ws.Range("L122").Value = 5158.799999999999
$ws.Range("M122").Value = -1019.3446
$ws.Range("N122").Value = -10058.8

# Sheet LTW, row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1572.4615
$ws.Range("I46").Value = 1984.5714
$ws.Range("J46").Value = 1091.6666
$ws.Range("K46").Value = 1984.5714
$ws.Range("L46").Value = 1091.6666
$ws.Range("M46").Value = -1796.5714
$ws.Range("N46").Value = -1467.6666

# Sheet LTW, row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2500
$ws.Range("I93").Value = 2600
$ws.Range("J93").Value = 2000
$ws.Range("K93").Value = 2600
$ws.Range("L93").Value = 2000
$ws.Range("M93").Value = -1352
$ws.Range("N93").Value = -4496

# Sheet WVR, row 47
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 14534.5
$ws.Range("J47").Value = 14534.5
$ws.Range("L47").Value = 14534.5
$ws.Range("N47").Value = -15678.5

# Sheet WVR, row 52
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 25000
$ws.Range("I52").Value = 9000
$ws.Range("J52").Value = 41000
$ws.Range("K52").Value = 9000
$ws.Range("L52").Value = 41000
$ws.Range("M52").Value = -8774
$ws.Range("N52").Value = -41452

# Sheet WVR, row 54
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 13086.5
$ws.Range("I54").Value = 1070
$ws.Range("J54").Value = 15489.8
$ws.Range("K54").Value = 1070
$ws.Range("L54").Value = 15489.8
$ws.Range("M54").Value = -550
$ws.Range("N54").Value = -16529.8

# Sheet WVR, row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 100000510
$ws.Range("I113").Value = 275
$ws.Range("J113").Value = 166667330
$ws.Range("K113").Value = 825
$ws.Range("L113").Value = 500001990
$ws.Range("M113").Value = 1345
$ws.Range("N113").Value = -500006330
